# Apply "repull data, push all data, mean calculation" update:
# Column F (dSF) values were repulled/recalculated for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = -6
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -4
$ws.Range("F12").Value = -7
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -4
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = -6
$ws.Range("F20").Value = 4
